# Add Venezuela to the country list, inserted as the new row 56
# (pushing the existing South Africa row down to row 57), matching
# the formatting style of similarly hand-added rows (e.g. Turkey, row 53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 56, shifting South Africa (and
# everything below it) down by one row.
$ws.Rows.Item(56).Insert()

# Copy the formatting (styles/number formats) only for columns A:F from
# row 53 (Turkey), which already uses "General" number formatting on every
# column - the same pattern used by this newly hand-entered row.
$ws.Range("A53:F53").Copy()
$ws.Range("A56:F56").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row with Venezuela's data.
$ws.Range("A56").Value = "Venezuela"
$ws.Range("B56").Value = "VEN"
$ws.Range("C56").Value = "VE"
$ws.Range("D56").Value = 10
$ws.Range("E56").Value = -66.9036
$ws.Range("F56").Value = 10.4806

# Reset the view: scroll back to the top and select cell A1, clearing the
# previous scroll position / selection that was left over from editing.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
